# material.xlsx maintenance edit:
#  - Capitalize the three sheet names (warehouse/region/transportation ->
#    Warehouse/Region/Transportation).
#  - Move the active tab / selection from the Transportation sheet to the
#    Region sheet, and update each sheet's remembered selection.
#  - Re-style the Anchor/LTE label columns (A, B, F, G) on the
#    Transportation sheet rows 2-33 to match the bordered "table edge" look
#    already used by column C (reuse existing cell style, don't invent a
#    new one).

$wb = $excel.ActiveWorkbook

# --- Rename sheets (case fix) --------------------------------------------
$wb.Worksheets.Item("warehouse").Name = "Warehouse"
$wb.Worksheets.Item("region").Name = "Region"
$wb.Worksheets.Item("transportation").Name = "Transportation"

$wsWarehouse = $wb.Worksheets.Item("Warehouse")
$wsRegion = $wb.Worksheets.Item("Region")
$wsTransportation = $wb.Worksheets.Item("Transportation")

# --- Re-style columns A, B, F, G for rows 2-33 on Transportation ---------
# Column C already carries the target "bordered" cell style - copy its
# format (not its value) onto the other label columns so the existing
# style index is reused instead of a new one being minted.
$styleSource = $wsTransportation.Range("C2")
$styleSource.Copy()
$wsTransportation.Range("A2:B33").PasteSpecial(-4122)   # xlPasteFormats
$styleSource.Copy()
$wsTransportation.Range("F2:G33").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update selections on each sheet -------------------------------------
# Leave the Warehouse sheet's own remembered selection untouched (it is not
# part of this edit).

# Transportation: update the remembered selection, but don't leave the
# sheet active - Region ends up as the active tab.
$wsTransportation.Range("J31").Select()

# Region becomes the active sheet with a new selection.
$wsRegion.Activate()
$wsRegion.Range("D54").Select()
